{"js": "// Apply the \"Use case specification - Place Order\" edit:\n//   1. Add a first-line indent of 720 twips (36 pt) to 17 specific body\n//      paragraphs (keeping their existing w:left=\"0\").\n//   2. In three of those paragraphs, Word's grammar checker additionally\n//      wrapped a subject/verb-agreement word (or phrase) in\n//      <w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>\n//      markers, which also splits the run that contained it.\n\n// Paragraphs that only need the first-line indent.\nconst plainIndentTargets = [\n  \"UC002\",\n  \"3.1 Customer\",\n  \"There is an active network connection to the Internet.\",\n  \"1. The customer request to place an order on the view cart screen\",\n  \"2. The AIMS software checks the availability of products in the cart\",\n  \"3. The AIMS software displays the form of delivery information\",\n  \"4. The customer enters and submits delivery information\",\n  \"6. The AIMS software calculates shipping fees\",\n  \"7. The AIMS software displays the invoice\",\n  \"9. The AIMS software calls UC \\u201cPay order\\u201d\",\n  \"10. The AIMS software saves order\",\n  \"11. The AIMS software makes the cart empty\",\n  \"12. The AIMS software displays the successful order notification\",\n  \"The logs have been updated accordingly\",\n];\n\n// Paragraphs that need both the indent AND a grammar-check run split.\n// Each entry gives the paragraph's current full text (to locate it) and the\n// OOXML body markup that should replace its contents (pPr is re-applied\n// afterwards through the regular paragraph-format API, so it only needs to\n// carry pStyle/jc here).\nconst splitTargets = [\n  {\n    text:\n      \"This use case describes the interaction between the AIMS software with the customer when the customer wishes to place order.\",\n    ooxmlParagraph:\n      '<w:p>' +\n      '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"left\"/></w:pPr>' +\n      '<w:r><w:rPr/><w:t xml:space=\"preserve\">This use case describes the interaction between the AIMS software with the customer when the customer wishes to place </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr/><w:t>order</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:rPr/><w:t>.</w:t></w:r>' +\n      '</w:p>',\n  },\n  {\n    text: \"5. The AIMS software check the validity of delivery information\",\n    ooxmlParagraph:\n      '<w:p>' +\n      '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"left\"/></w:pPr>' +\n      '<w:r><w:rPr/><w:t xml:space=\"preserve\">5. The AIMS software </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr/><w:t>check</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:rPr/><w:t xml:space=\"preserve\"> the validity of delivery information</w:t></w:r>' +\n      '</w:p>',\n  },\n  {\n    text: \"8. The customer confirms to place order\",\n    ooxmlParagraph:\n      '<w:p>' +\n      '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"left\"/></w:pPr>' +\n      '<w:r><w:rPr/><w:t xml:space=\"preserve\">8. The customer confirms </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr/><w:t>to place</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:rPr/><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr/><w:t>order</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '</w:p>',\n  },\n];\n\nfunction wrapOoxml(paragraphXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    paragraphXml +\n    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nconst body = context.document.body;\n\n// --- Pass 1: handle the three paragraphs that need a run split first, while\n// their text still matches the *original* wording (insertOoxml replaces the\n// whole paragraph, including its pPr, so we restore the indent in pass 2). ---\nfor (const { text, ooxmlParagraph } of splitTargets) {\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n\n  const match = paras.items.find((p) => p.text === text);\n  if (!match) {\n    throw new Error(\"Paragraph not found: \" + text);\n  }\n  match.insertOoxml(wrapOoxml(ooxmlParagraph), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Pass 2: apply the first-line indent (w:left stays 0) to every target\n// paragraph, re-locating each by its (possibly just-rewritten) text. ---\nconst allTargetTexts = new Set([\n  ...plainIndentTargets,\n  ...splitTargets.map((t) => t.text),\n]);\n\nconst finalParas = body.paragraphs;\nfinalParas.load(\"text\");\nawait context.sync();\n\nfor (const p of finalParas.items) {\n  if (allTargetTexts.has(p.text)) {\n    p.leftIndent = 0;\n    p.firstLineIndent = 36; // 720 twips == 36 points\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the \"Use case specification - Place Order\" edit:\n#   1. Add a first-line indent of 720 twips (36 pt) to 17 specific body\n#      paragraphs (keeping their existing left indent of 0).\n#   2. In three of those paragraphs, Word's grammar checker additionally\n#      wrapped a subject/verb-agreement word (or phrase) in\n#      <w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>\n#      markers, which also splits the run that contained it.\n\n$d = $word.ActiveDocument\n\n# Paragraphs that need only the first-line indent (matched by exact text).\n$plainIndentTargets = @(\n    \"UC002\",\n    \"3.1 Customer\",\n    \"There is an active network connection to the Internet.\",\n    \"1. The customer request to place an order on the view cart screen\",\n    \"2. The AIMS software checks the availability of products in the cart\",\n    \"3. The AIMS software displays the form of delivery information\",\n    \"4. The customer enters and submits delivery information\",\n    \"6. The AIMS software calculates shipping fees\",\n    \"7. The AIMS software displays the invoice\",\n    (\"9. The AIMS software calls UC \" + [char]0x201C + \"Pay order\" + [char]0x201D),\n    \"10. The AIMS software saves order\",\n    \"11. The AIMS software makes the cart empty\",\n    \"12. The AIMS software displays the successful order notification\",\n    \"The logs have been updated accordingly\"\n)\n\n# Paragraphs that need both the indent AND a grammar-check run split. Maps\n# the paragraph's current exact text to the OOXML that should replace its\n# content (pPr here only needs pStyle/jc - the indent is re-applied via the\n# ParagraphFormat API afterwards so w:left=\"0\" survives).\n$splitTargets = [ordered]@{\n    \"This use case describes the interaction between the AIMS software with the customer when the customer wishes to place order.\" =\n        '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"left\"/></w:pPr>' +\n        '<w:r><w:rPr/><w:t xml:space=\"preserve\">This use case describes the interaction between the AIMS software with the customer when the customer wishes to place </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:rPr/><w:t>order</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:rPr/><w:t>.</w:t></w:r>' +\n        '</w:p>';\n    \"5. The AIMS software check the validity of delivery information\" =\n        '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"left\"/></w:pPr>' +\n        '<w:r><w:rPr/><w:t xml:space=\"preserve\">5. The AIMS software </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:rPr/><w:t>check</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:rPr/><w:t xml:space=\"preserve\"> the validity of delivery information</w:t></w:r>' +\n        '</w:p>';\n    \"8. The customer confirms to place order\" =\n        '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:pPr><w:pStyle w:val=\"Normal\"/><w:jc w:val=\"left\"/></w:pPr>' +\n        '<w:r><w:rPr/><w:t xml:space=\"preserve\">8. The customer confirms </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:rPr/><w:t>to place</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:rPr/><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:rPr/><w:t>order</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '</w:p>'\n}\n\n# Union of every paragraph's \"final\" text, used in the indent-fix-up pass.\n$allTargetTexts = @{}\nforeach ($t in $plainIndentTargets) { $allTargetTexts[$t] = $true }\nforeach ($k in $splitTargets.Keys) { $allTargetTexts[$k] = $true }\n\n# --- Pass 1: run-split the three grammar-flagged paragraphs first, while\n# their text still matches the original wording. InsertXML replaces the\n# paragraph's content in place (paragraph count/index is unaffected). ---\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]0x0D, [char]0x07)\n    if ($splitTargets.Contains($t)) {\n        [void]$p.Range.InsertXML($splitTargets[$t])\n    }\n}\n\n# --- Pass 2: apply the first-line indent (left indent stays 0) to every\n# target paragraph. ---\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]0x0D, [char]0x07)\n    if ($allTargetTexts.ContainsKey($t)) {\n        $p.Format.LeftIndent = 0\n        $p.Format.FirstLineIndent = 36\n    }\n}\n"}
